# feat: add 2022-Q1 data
#
# The workbook tracks quarterly fund-holdings snapshots, one sheet per
# quarter, plus a trailing "总计" (totals) summary sheet. Each time a new
# quarter is published, the old "总计" sheet is renamed to the new quarter
# and filled with that quarter's fund-holdings detail, and a brand new
# "总计" sheet is appended with the updated summary history.

$wb = $excel.ActiveWorkbook

# A sheet that already carries the workbook's standard header/index style
# (bold font + thin box border, centered) - used below as a formatting
# template so we reuse the existing style entries instead of minting new
# ones via Font/Borders property writes.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# ------------------------------------------------------------------
# 1) Turn the existing "总计" sheet into the new "2022-Q1" detail sheet,
#    then append a fresh "总计" sheet right after it.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# ------------------------------------------------------------------
# 2) Populate "2022-Q1" with the fund-holdings detail table.
# ------------------------------------------------------------------

# Clear out whatever the sheet used to hold (old "总计" summary rows).
$q1.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

$data = @(
    ,@("163409", "兴全绿色投资混合(LOF)", "77.81", "89.66", "7.15", "5.5634", 2)
    ,@("110029", "易方达科讯混合", "36.09", "91.34", "8.82", "3.1831", 1)
    ,@("006533", "易方达科融混合", "32.14", "89.61", "8.56", "2.7512", 1)
    ,@("163412", "兴全轻资产混合(LOF)", "69.20", "95.26", "3.61", "2.4981", 8)
    ,@("110013", "易方达科翔混合", "49.44", "87.59", "4.04", "1.9974", 4)
    ,@("009341", "易方达均衡成长股票", "64.76", "82.97", "2.90", "1.8780", 7)
    ,@("003293", "易方达科瑞灵活配置混合", "34.67", "78.17", "3.27", "1.1337", 5)
    ,@("001513", "易方达信息产业混合", "32.50", "92.37", "3.35", "1.0888", 3)
    ,@("900003", "中信证券臻选价值成长混合A", "48.09", "86.76", "2.03", "0.9762", 8)
    ,@("110001", "易方达平稳增长混合", "33.39", "60.98", "2.62", "0.8748", 9)
    ,@("010013", "易方达信息行业精选股票", "25.38", "88.64", "3.05", "0.7741", 6)
    ,@("340006", "兴全全球视野股票", "21.42", "83.30", "3.00", "0.6426", 10)
    ,@("110012", "易方达科汇灵活配置混合", "15.73", "75.64", "3.04", "0.4782", 5)
    ,@("011649", "易方达逆向投资混合A", "7.49", "85.02", "3.23", "0.2419", 7)
    ,@("160529", "博时创业板两年定期开放混合", "7.92", "82.61", "2.84", "0.2249", 9)
    ,@("010389", "易方达科益混合A", "7.10", "92.94", "3.16", "0.2244", 10)
    ,@("001521", "国寿安保成长优选股票", "4.24", "87.83", "5.18", "0.2196", 3)
    ,@("900079", "中信证券臻选价值成长混合C", "10.59", "86.76", "2.03", "0.2150", 8)
    ,@("001990", "中欧数据挖掘多因子灵活配置混合A", "18.03", "84.18", "0.76", "0.1370", 3)
    ,@("014135", "中欧金安量化混合A", "9.43", "67.44", "0.69", "0.0651", 3)
    ,@("011650", "易方达逆向投资混合C", "1.96", "85.02", "3.23", "0.0633", 7)
    ,@("159804", "国寿安保国证创业板中盘精选88ETF", "2.10", "98.79", "1.94", "0.0407", 8)
    ,@("004234", "中欧数据挖掘多因子灵活配置混合C", "5.06", "84.18", "0.76", "0.0385", 3)
    ,@("008082", "国寿安保研究精选混合A", "0.52", "91.60", "5.60", "0.0291", 4)
    ,@("010390", "易方达科益混合C", "0.29", "92.94", "3.16", "0.0092", 10)
    ,@("014136", "中欧金安量化混合C", "1.28", "67.44", "0.69", "0.0088", 3)
    ,@("008083", "国寿安保研究精选混合C", "0.15", "91.60", "5.60", "0.0084", 4)
)

$lastRow = $data.Count + 1

# Header row (B1:H1), bold+bordered like every other quarter sheet.
for ($c = 2; $c -le 8; $c++) {
    $q1.Cells.Item(1, $c).Value = $headers[$c - 2]
}

# Column A holds the 0-based row index, bold+bordered, numeric.
for ($i = 0; $i -lt $data.Count; $i++) {
    $q1.Cells.Item($i + 2, 1).Value = $i
}

# Columns B-G are stored as plain text (fund code/name/percentages), column H
# (rank) is numeric - matching the source data's own typing.
$q1.Range("B2:G$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 0; $c -lt 6; $c++) {
        $q1.Cells.Item($r, $c + 2).Value = $row[$c]
    }
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# Re-apply the workbook's existing header/index style (bold font + thin box
# border, centered) by copying it from a known styled cell instead of
# re-building fonts/borders (which would mint new style entries).
$styleSrc.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$q1.Range("A2:A$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q1.Range("A1").Select()

# ------------------------------------------------------------------
# 3) Populate the new "总计" sheet: the same summary history, with a new
#    leading row for 2022-Q1 and everything else shifted down by one.
# ------------------------------------------------------------------
$total.Cells.Clear()

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$summary = @(
    ,@("2022-Q1", 27, 25.37)
    ,@("2021-Q4", 19, 24.79)
    ,@("2021-Q3", 30, 26.43)
    ,@("2021-Q2", 33, 18.7)
    ,@("2021-Q1", 34, 20.18)
    ,@("2020-Q4", 6, 11.03)
)

for ($i = 0; $i -lt $summary.Count; $i++) {
    $row = $summary[$i]
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

$styleSrc.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$total.Range("A2:A$($summary.Count + 1)").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A1").Select()
